$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.008.38"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.102.04"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.16"
$ws.Range("E5").Value = "  +3.69%  "
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5149"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4429"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.28"
$ws.Range("E9").Value = "  -4.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08966"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.169"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.099.62"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.206"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.734"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "98.87"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("E19").Value = "  +7.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06668"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.214"
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.120.19"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.62"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.338"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.348.75"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.553"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.78"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.21"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.639"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.208"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.122"
$ws.Range("E36").Value = "  +4.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.13"
$ws.Range("E37").Value = "  -4.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02567"
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06781"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2281"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.349"
$ws.Range("E41").Value = "  +7.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.48"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6794"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.11"
$ws.Range("E44").Value = "  -4.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6379"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000365"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.640"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.20"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07218"
$ws.Range("E51").Value = "  -0.21%  "
